# RPA datasets push 2023-12-01
# Insert a new IPO row ("포스뱅크") above the current top data row,
# then drop the obsolete "키움스팩9호" row that the insert pushed down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by one row, starting at row 2 (right below the header).
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new listing's data.
$ws.Range("A2").Value = "포스뱅크"
$ws.Range("B2").Value = "2024.01.05~01.11"
$ws.Range("C2").Value = "13,000~15,000"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 19500
$ws.Range("F2").Value = "하나증권"

# The insert pushed the old "키움스팩9호" row from row 3 to row 4 - remove it.
$ws.Rows("4:4").Delete()
